# ============================================================
# Refresh crypto market data in the "Live Data" sheet and
# recompute the summary figures in the "Analysis" sheet.
# (Add files via upload)
# ============================================================

# Helper: assign a value as literal text, guarding against Excel
# auto-converting strings that look like numbers/currency/dates
# into numeric cells, then restore the default (no) style so the
# cell format stays exactly as it was before.
function Set-TextValue {
    param($Cell, $Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$wsLive = $wb.Worksheets.Item("Live Data")
$wsAnalysis = $wb.Worksheets.Item("Analysis")

# ---- Live Data sheet: rows 2-51 (Name, Symbol, Price, MarketCap, Volume, Change%) ----

# Row 2: Bitcoin (BTC)
$wsLive.Cells.Item(2, 3).Value = 84691.0
$wsLive.Cells.Item(2, 4).Value = 1668305607915.0
$wsLive.Cells.Item(2, 5).Value = 68615162383.0
$wsLive.Cells.Item(2, 6).Value = -6.93496

# Row 3: Ethereum (ETH)
$wsLive.Cells.Item(3, 3).Value = 2126.87
$wsLive.Cells.Item(3, 4).Value = 254574919263.0
$wsLive.Cells.Item(3, 5).Value = 36808605728.0
$wsLive.Cells.Item(3, 6).Value = -8.34483

# Row 4: Tether (USDT)
$wsLive.Cells.Item(4, 3).Value = 0.999333
$wsLive.Cells.Item(4, 4).Value = 142414725432.0
$wsLive.Cells.Item(4, 5).Value = 90993221888.0
$wsLive.Cells.Item(4, 6).Value = -0.03518

# Row 5: XRP (XRP)
$wsLive.Cells.Item(5, 3).Value = 2.44
$wsLive.Cells.Item(5, 4).Value = 140234281340.0
$wsLive.Cells.Item(5, 5).Value = 13571971959.0
$wsLive.Cells.Item(5, 6).Value = -7.69101

# Row 6: BNB (BNB)
$wsLive.Cells.Item(6, 3).Value = 577.21
$wsLive.Cells.Item(6, 4).Value = 83736229850.0
$wsLive.Cells.Item(6, 5).Value = 1171766441.0
$wsLive.Cells.Item(6, 6).Value = -2.95984

# Row 7: Solana (SOL)
$wsLive.Cells.Item(7, 3).Value = 142.28
$wsLive.Cells.Item(7, 4).Value = 71509536107.0
$wsLive.Cells.Item(7, 5).Value = 10107983974.0
$wsLive.Cells.Item(7, 6).Value = -11.25725

# Row 8: USDC (USDC)
$wsLive.Cells.Item(8, 3).Value = 0.999968
$wsLive.Cells.Item(8, 4).Value = 56328047113.0
$wsLive.Cells.Item(8, 5).Value = 10963551917.0
$wsLive.Cells.Item(8, 6).Value = 0.00884

# Row 9: Cardano (ADA)
$wsLive.Cells.Item(9, 3).Value = 0.859436
$wsLive.Cells.Item(9, 4).Value = 30705483779.0
$wsLive.Cells.Item(9, 5).Value = 5379900666.0
$wsLive.Cells.Item(9, 6).Value = -11.01596

# Row 10: Dogecoin (DOGE)
$wsLive.Cells.Item(10, 3).Value = 0.198158
$wsLive.Cells.Item(10, 4).Value = 29072112585.0
$wsLive.Cells.Item(10, 5).Value = 2888601783.0
$wsLive.Cells.Item(10, 6).Value = -8.94008

# Row 11: TRON (TRX)
$wsLive.Cells.Item(11, 3).Value = 0.237121
$wsLive.Cells.Item(11, 4).Value = 20305505082.0
$wsLive.Cells.Item(11, 5).Value = 1413304532.0
$wsLive.Cells.Item(11, 6).Value = -2.12085

# Row 12: Lido Staked Ether (STETH)
$wsLive.Cells.Item(12, 3).Value = 2119.77
$wsLive.Cells.Item(12, 4).Value = 19707926200.0
$wsLive.Cells.Item(12, 5).Value = 183294670.0
$wsLive.Cells.Item(12, 6).Value = -8.34596

# Row 13: Pi Network (PI)
Set-TextValue $wsLive.Cells.Item(13, 1) "Pi Network"
Set-TextValue $wsLive.Cells.Item(13, 2) "PI"
$wsLive.Cells.Item(13, 3).Value = 1.78
$wsLive.Cells.Item(13, 4).Value = 12260513497.0
$wsLive.Cells.Item(13, 5).Value = 719040012.0
$wsLive.Cells.Item(13, 6).Value = 3.50195

# Row 14: Wrapped Bitcoin (WBTC)
Set-TextValue $wsLive.Cells.Item(14, 1) "Wrapped Bitcoin"
Set-TextValue $wsLive.Cells.Item(14, 2) "WBTC"
$wsLive.Cells.Item(14, 3).Value = 84601.0
$wsLive.Cells.Item(14, 4).Value = 10925312028.0
$wsLive.Cells.Item(14, 5).Value = 821737492.0
$wsLive.Cells.Item(14, 6).Value = -6.8824

# Row 15: Hedera (HBAR)
Set-TextValue $wsLive.Cells.Item(15, 1) "Hedera"
Set-TextValue $wsLive.Cells.Item(15, 2) "HBAR"
$wsLive.Cells.Item(15, 3).Value = 0.237376
$wsLive.Cells.Item(15, 4).Value = 9915148974.0
$wsLive.Cells.Item(15, 5).Value = 989785725.0
$wsLive.Cells.Item(15, 6).Value = -5.33909

# Row 16: LEO Token (LEO)
Set-TextValue $wsLive.Cells.Item(16, 1) "LEO Token"
Set-TextValue $wsLive.Cells.Item(16, 2) "LEO"
$wsLive.Cells.Item(16, 3).Value = 10.07
$wsLive.Cells.Item(16, 4).Value = 9305895793.0
$wsLive.Cells.Item(16, 5).Value = 8609395.0
$wsLive.Cells.Item(16, 6).Value = 1.7533

# Row 17: Wrapped stETH (WSTETH)
Set-TextValue $wsLive.Cells.Item(17, 1) "Wrapped stETH"
Set-TextValue $wsLive.Cells.Item(17, 2) "WSTETH"
$wsLive.Cells.Item(17, 3).Value = 2563.26
$wsLive.Cells.Item(17, 4).Value = 9012985962.0
$wsLive.Cells.Item(17, 5).Value = 77443074.0
$wsLive.Cells.Item(17, 6).Value = -7.39994

# Row 18: Chainlink (LINK)
Set-TextValue $wsLive.Cells.Item(18, 1) "Chainlink"
Set-TextValue $wsLive.Cells.Item(18, 2) "LINK"
$wsLive.Cells.Item(18, 3).Value = 14.07
$wsLive.Cells.Item(18, 4).Value = 8912052672.0
$wsLive.Cells.Item(18, 5).Value = 1085512749.0
$wsLive.Cells.Item(18, 6).Value = -14.91149

# Row 19: Stellar (XLM)
Set-TextValue $wsLive.Cells.Item(19, 1) "Stellar"
Set-TextValue $wsLive.Cells.Item(19, 2) "XLM"
$wsLive.Cells.Item(19, 3).Value = 0.289289
$wsLive.Cells.Item(19, 4).Value = 8806615176.0
$wsLive.Cells.Item(19, 5).Value = 422558629.0
$wsLive.Cells.Item(19, 6).Value = -11.25966

# Row 20: Avalanche (AVAX)
Set-TextValue $wsLive.Cells.Item(20, 1) "Avalanche"
Set-TextValue $wsLive.Cells.Item(20, 2) "AVAX"
$wsLive.Cells.Item(20, 3).Value = 20.4
$wsLive.Cells.Item(20, 4).Value = 8391112095.0
$wsLive.Cells.Item(20, 5).Value = 822337550.0
$wsLive.Cells.Item(20, 6).Value = -12.83792

# Row 21: USDS (USDS)
Set-TextValue $wsLive.Cells.Item(21, 1) "USDS"
Set-TextValue $wsLive.Cells.Item(21, 2) "USDS"
$wsLive.Cells.Item(21, 3).Value = 1.0
$wsLive.Cells.Item(21, 4).Value = 7945048807.0
$wsLive.Cells.Item(21, 5).Value = 61652474.0
$wsLive.Cells.Item(21, 6).Value = 0.0408

# Row 22: Litecoin (LTC)
Set-TextValue $wsLive.Cells.Item(22, 1) "Litecoin"
Set-TextValue $wsLive.Cells.Item(22, 2) "LTC"
$wsLive.Cells.Item(22, 3).Value = 104.33
$wsLive.Cells.Item(22, 4).Value = 7812184362.0
$wsLive.Cells.Item(22, 5).Value = 1574325261.0
$wsLive.Cells.Item(22, 6).Value = -9.56579

# Row 23: Toncoin (TON)
Set-TextValue $wsLive.Cells.Item(23, 1) "Toncoin"
Set-TextValue $wsLive.Cells.Item(23, 2) "TON"
$wsLive.Cells.Item(23, 3).Value = 3.12
$wsLive.Cells.Item(23, 4).Value = 7678972297.0
$wsLive.Cells.Item(23, 5).Value = 320270192.0
$wsLive.Cells.Item(23, 6).Value = -5.52245

# Row 24: Sui (SUI)
Set-TextValue $wsLive.Cells.Item(24, 1) "Sui"
Set-TextValue $wsLive.Cells.Item(24, 2) "SUI"
$wsLive.Cells.Item(24, 3).Value = 2.45
$wsLive.Cells.Item(24, 4).Value = 7670848221.0
$wsLive.Cells.Item(24, 5).Value = 1546831486.0
$wsLive.Cells.Item(24, 6).Value = -16.24126

# Row 25: Shiba Inu (SHIB)
Set-TextValue $wsLive.Cells.Item(25, 1) "Shiba Inu"
Set-TextValue $wsLive.Cells.Item(25, 2) "SHIB"
$wsLive.Cells.Item(25, 3).Value = 0.00001305
$wsLive.Cells.Item(25, 4).Value = 7633476365.0
$wsLive.Cells.Item(25, 5).Value = 514222609.0
$wsLive.Cells.Item(25, 6).Value = -6.26859

# Row 26: MANTRA (OM)
Set-TextValue $wsLive.Cells.Item(26, 1) "MANTRA"
Set-TextValue $wsLive.Cells.Item(26, 2) "OM"
$wsLive.Cells.Item(26, 3).Value = 7.06
$wsLive.Cells.Item(26, 4).Value = 6820733178.0
$wsLive.Cells.Item(26, 5).Value = 203414338.0
$wsLive.Cells.Item(26, 6).Value = -3.01929

# Row 27: Polkadot (DOT)
Set-TextValue $wsLive.Cells.Item(27, 1) "Polkadot"
Set-TextValue $wsLive.Cells.Item(27, 2) "DOT"
$wsLive.Cells.Item(27, 3).Value = 4.3
$wsLive.Cells.Item(27, 4).Value = 6484332788.0
$wsLive.Cells.Item(27, 5).Value = 410581788.0
$wsLive.Cells.Item(27, 6).Value = -12.2823

# Row 28: Bitcoin Cash (BCH)
Set-TextValue $wsLive.Cells.Item(28, 1) "Bitcoin Cash"
Set-TextValue $wsLive.Cells.Item(28, 2) "BCH"
$wsLive.Cells.Item(28, 3).Value = 316.2
$wsLive.Cells.Item(28, 4).Value = 6213730825.0
$wsLive.Cells.Item(28, 5).Value = 491830498.0
$wsLive.Cells.Item(28, 6).Value = -6.09046

# Row 29: WETH (WETH)
Set-TextValue $wsLive.Cells.Item(29, 1) "WETH"
Set-TextValue $wsLive.Cells.Item(29, 2) "WETH"
$wsLive.Cells.Item(29, 3).Value = 2138.05
$wsLive.Cells.Item(29, 4).Value = 6079376797.0
$wsLive.Cells.Item(29, 5).Value = 1901220443.0
$wsLive.Cells.Item(29, 6).Value = -7.79625

# Row 30: Ethena USDe (USDE)
$wsLive.Cells.Item(30, 3).Value = 0.998899
$wsLive.Cells.Item(30, 4).Value = 5448284065.0
$wsLive.Cells.Item(30, 5).Value = 104373459.0
$wsLive.Cells.Item(30, 6).Value = 0.0468

# Row 31: Hyperliquid (HYPE)
Set-TextValue $wsLive.Cells.Item(31, 1) "Hyperliquid"
Set-TextValue $wsLive.Cells.Item(31, 2) "HYPE"
$wsLive.Cells.Item(31, 3).Value = 16.42
$wsLive.Cells.Item(31, 4).Value = 5434103062.0
$wsLive.Cells.Item(31, 5).Value = 368941150.0
$wsLive.Cells.Item(31, 6).Value = -15.62228

# Row 32: Bitget Token (BGB)
Set-TextValue $wsLive.Cells.Item(32, 1) "Bitget Token"
Set-TextValue $wsLive.Cells.Item(32, 2) "BGB"
$wsLive.Cells.Item(32, 3).Value = 4.17
$wsLive.Cells.Item(32, 4).Value = 4921508738.0
$wsLive.Cells.Item(32, 5).Value = 353213116.0
$wsLive.Cells.Item(32, 6).Value = -7.74536

# Row 33: Wrapped eETH (WEETH)
Set-TextValue $wsLive.Cells.Item(33, 1) "Wrapped eETH"
Set-TextValue $wsLive.Cells.Item(33, 2) "WEETH"
$wsLive.Cells.Item(33, 3).Value = 2263.01
$wsLive.Cells.Item(33, 4).Value = 4334128454.0
$wsLive.Cells.Item(33, 5).Value = 30383710.0
$wsLive.Cells.Item(33, 6).Value = -7.87038

# Row 34: WhiteBIT Coin (WBT)
Set-TextValue $wsLive.Cells.Item(34, 1) "WhiteBIT Coin"
Set-TextValue $wsLive.Cells.Item(34, 2) "WBT"
$wsLive.Cells.Item(34, 3).Value = 29.93
$wsLive.Cells.Item(34, 4).Value = 4281665721.0
$wsLive.Cells.Item(34, 5).Value = 141534155.0
$wsLive.Cells.Item(34, 6).Value = -1.56903

# Row 35: Uniswap (UNI)
Set-TextValue $wsLive.Cells.Item(35, 1) "Uniswap"
Set-TextValue $wsLive.Cells.Item(35, 2) "UNI"
$wsLive.Cells.Item(35, 3).Value = 6.92
$wsLive.Cells.Item(35, 4).Value = 4126709295.0
$wsLive.Cells.Item(35, 5).Value = 365903364.0
$wsLive.Cells.Item(35, 6).Value = -11.37376

# Row 36: Monero (XMR)
Set-TextValue $wsLive.Cells.Item(36, 1) "Monero"
Set-TextValue $wsLive.Cells.Item(36, 2) "XMR"
$wsLive.Cells.Item(36, 3).Value = 216.82
$wsLive.Cells.Item(36, 4).Value = 3964995023.0
$wsLive.Cells.Item(36, 5).Value = 65581423.0
$wsLive.Cells.Item(36, 6).Value = -5.42726

# Row 37: NEAR Protocol (NEAR)
Set-TextValue $wsLive.Cells.Item(37, 1) "NEAR Protocol"
Set-TextValue $wsLive.Cells.Item(37, 2) "NEAR"
$wsLive.Cells.Item(37, 3).Value = 2.81
$wsLive.Cells.Item(37, 4).Value = 3317122524.0
$wsLive.Cells.Item(37, 5).Value = 412998239.0
$wsLive.Cells.Item(37, 6).Value = -13.97688

# Row 38: Dai (DAI)
Set-TextValue $wsLive.Cells.Item(38, 1) "Dai"
Set-TextValue $wsLive.Cells.Item(38, 2) "DAI"
$wsLive.Cells.Item(38, 3).Value = 0.999606
$wsLive.Cells.Item(38, 4).Value = 3277977105.0
$wsLive.Cells.Item(38, 5).Value = 224501301.0
$wsLive.Cells.Item(38, 6).Value = 0.00693

# Row 39: Aptos (APT)
Set-TextValue $wsLive.Cells.Item(39, 1) "Aptos"
Set-TextValue $wsLive.Cells.Item(39, 2) "APT"
$wsLive.Cells.Item(39, 3).Value = 5.48
$wsLive.Cells.Item(39, 4).Value = 3205964591.0
$wsLive.Cells.Item(39, 5).Value = 389307904.0
$wsLive.Cells.Item(39, 6).Value = -12.94364

# Row 40: sUSDS (SUSDS)
Set-TextValue $wsLive.Cells.Item(40, 1) "sUSDS"
Set-TextValue $wsLive.Cells.Item(40, 2) "SUSDS"
$wsLive.Cells.Item(40, 3).Value = 1.042
$wsLive.Cells.Item(40, 4).Value = 3016171614.0
$wsLive.Cells.Item(40, 5).Value = 1759153.0
$wsLive.Cells.Item(40, 6).Value = 0.11259

# Row 41: Ondo (ONDO)
Set-TextValue $wsLive.Cells.Item(41, 1) "Ondo"
Set-TextValue $wsLive.Cells.Item(41, 2) "ONDO"
$wsLive.Cells.Item(41, 3).Value = 0.939233
$wsLive.Cells.Item(41, 4).Value = 2939621432.0
$wsLive.Cells.Item(41, 5).Value = 483909855.0
$wsLive.Cells.Item(41, 6).Value = -17.30967

# Row 42: Pepe (PEPE)
Set-TextValue $wsLive.Cells.Item(42, 1) "Pepe"
Set-TextValue $wsLive.Cells.Item(42, 2) "PEPE"
$wsLive.Cells.Item(42, 3).Value = 0.00000709
$wsLive.Cells.Item(42, 4).Value = 2937553429.0
$wsLive.Cells.Item(42, 5).Value = 783693527.0
$wsLive.Cells.Item(42, 6).Value = -12.95279

# Row 43: Aave (AAVE)
Set-TextValue $wsLive.Cells.Item(43, 1) "Aave"
Set-TextValue $wsLive.Cells.Item(43, 2) "AAVE"
$wsLive.Cells.Item(43, 3).Value = 195.16
$wsLive.Cells.Item(43, 4).Value = 2914938754.0
$wsLive.Cells.Item(43, 5).Value = 622343859.0
$wsLive.Cells.Item(43, 6).Value = -5.12873

# Row 44: Internet Computer (ICP)
Set-TextValue $wsLive.Cells.Item(44, 1) "Internet Computer"
Set-TextValue $wsLive.Cells.Item(44, 2) "ICP"
$wsLive.Cells.Item(44, 3).Value = 6.06
$wsLive.Cells.Item(44, 4).Value = 2891401150.0
$wsLive.Cells.Item(44, 5).Value = 128465197.0
$wsLive.Cells.Item(44, 6).Value = -10.67337

# Row 45: Ethereum Classic (ETC)
Set-TextValue $wsLive.Cells.Item(45, 1) "Ethereum Classic"
Set-TextValue $wsLive.Cells.Item(45, 2) "ETC"
$wsLive.Cells.Item(45, 3).Value = 18.77
$wsLive.Cells.Item(45, 4).Value = 2814881389.0
$wsLive.Cells.Item(45, 5).Value = 196056069.0
$wsLive.Cells.Item(45, 6).Value = -5.52695

# Row 46: Gate (GT)
Set-TextValue $wsLive.Cells.Item(46, 1) "Gate"
Set-TextValue $wsLive.Cells.Item(46, 2) "GT"
$wsLive.Cells.Item(46, 3).Value = 20.19
$wsLive.Cells.Item(46, 4).Value = 2540332856.0
$wsLive.Cells.Item(46, 5).Value = 24405956.0
$wsLive.Cells.Item(46, 6).Value = -6.6648

# Row 47: Official Trump (TRUMP)
$wsLive.Cells.Item(47, 3).Value = 12.59
$wsLive.Cells.Item(47, 4).Value = 2513767633.0
$wsLive.Cells.Item(47, 5).Value = 2122514701.0
$wsLive.Cells.Item(47, 6).Value = -15.46241

# Row 48: OKB (OKB)
Set-TextValue $wsLive.Cells.Item(48, 1) "OKB"
Set-TextValue $wsLive.Cells.Item(48, 2) "OKB"
$wsLive.Cells.Item(48, 3).Value = 41.8
$wsLive.Cells.Item(48, 4).Value = 2508178817.0
$wsLive.Cells.Item(48, 5).Value = 29641948.0
$wsLive.Cells.Item(48, 6).Value = -8.57925

# Row 49: Coinbase Wrapped BTC (CBBTC)
Set-TextValue $wsLive.Cells.Item(49, 1) "Coinbase Wrapped BTC"
Set-TextValue $wsLive.Cells.Item(49, 2) "CBBTC"
$wsLive.Cells.Item(49, 3).Value = 84426.0
$wsLive.Cells.Item(49, 4).Value = 2425794323.0
$wsLive.Cells.Item(49, 5).Value = 526518365.0
$wsLive.Cells.Item(49, 6).Value = -7.31215

# Row 50: Mantle (MNT)
Set-TextValue $wsLive.Cells.Item(50, 1) "Mantle"
Set-TextValue $wsLive.Cells.Item(50, 2) "MNT"
$wsLive.Cells.Item(50, 3).Value = 0.727613
$wsLive.Cells.Item(50, 4).Value = 2422127128.0
$wsLive.Cells.Item(50, 5).Value = 116903037.0
$wsLive.Cells.Item(50, 6).Value = -4.04246

# Row 51: Bittensor (TAO)
Set-TextValue $wsLive.Cells.Item(51, 1) "Bittensor"
Set-TextValue $wsLive.Cells.Item(51, 2) "TAO"
$wsLive.Cells.Item(51, 3).Value = 277.62
$wsLive.Cells.Item(51, 4).Value = 2316769533.0
$wsLive.Cells.Item(51, 5).Value = 222940774.0
$wsLive.Cells.Item(51, 6).Value = -14.15931

# ---- Analysis sheet: refresh summary metrics ----

Set-TextValue $wsAnalysis.Cells.Item(2, 2) "2025-03-04 20:19:17"   # Last Updated
Set-TextValue $wsAnalysis.Cells.Item(4, 2) "`$5339.98"   # Average Price (USD)
Set-TextValue $wsAnalysis.Cells.Item(5, 2) "Pi Network (3.50%)"   # Highest 24h Change
Set-TextValue $wsAnalysis.Cells.Item(6, 2) "Ondo (-17.31%)"   # Lowest 24h Change
